# The "KIQ-8710" vehicle (rows 58-59) is being removed from the
# informativo / cleanup listing. Deleting the two rows shifts every
# following row up by two, which matches the target state exactly
# (shared strings, styles and values all ride along with their rows).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A58:A59").EntireRow.Delete()

# Leave the selection where the author ended up after the edit.
$ws.Range("G65").Select()
